$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 1566
$ws1.Range("F10").Value = 3480
$ws1.Range("F14").Value = 25
$ws1.Range("F16").Value = 9
$ws1.Range("F17").Value = 1196
$ws1.Range("F23").Value = 1588

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F23").Value = 99
$ws2.Range("F36").Value = 418
$ws2.Range("F41").Value = 6

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F13").Value = 782
$ws3.Range("F14").Value = 171

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 782
$ws4.Range("F12").Value = 1566
$ws4.Range("F20").Value = 25
$ws4.Range("F24").Value = 1196
$ws4.Range("F37").Value = 99
$ws4.Range("F38").Value = 99
$ws4.Range("F41").Value = 1588
$ws4.Range("F46").Value = 418
